# Apply the "created family Ackermann to test new variable dossier_id" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at N:O (shifts old N..R -> P..T) -----------------
$ws.Range("N1:O1").EntireColumn.Insert()

# --- Column B width (widened to fit the new "Ackermann-Abegger" surname) --------
$ws.Columns.Item(2).ColumnWidth = 17.2

# --- New header cells ----------------------------------------------------------
$ws.Range("N1").Value = "EGID"
$ws.Range("O1").Value = "EWID"

# --- Row 2 / 3 tweaks (civil status values) ------------------------------------
$ws.Range("F2").Value = "Civil stat#B"
$ws.Range("F3").Value = "Civil stat#L"

# --- Row 2 new EGID/EWID values -------------------------------------------------
$ws.Range("N2").Value = 77777
$ws.Range("O2").Value = 1

# --- Row 3 new EGID/EWID values -------------------------------------------------
$ws.Range("N3").Value = 6666
$ws.Range("O3").Value = 11

# --- Row 4 new EGID/EWID values (must carry the row's special style) ------------
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("N4").Value = 9999999
$ws.Range("O4").Value = 999

# --- New row 9: Anton Ackermann --------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("A9").Value = 7560000000007
$ws.Range("B9").Value = "Ackermann"
$ws.Range("C9").Value = "Anton"
$ws.Range("D9").Value = 34940
$ws.Range("E9").Value = "Sex#M"
$ws.Range("F9").Value = "Civil stat#B"
$ws.Range("G9").Value = "Antragssteller"
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = "CHResidenc#B"
$ws.Range("J9").Value = "Ackerstrasse"
$ws.Range("K9").Value = 11
$ws.Range("L9").Value = 804500
$ws.Range("M9").Value = "Zürich"
$ws.Range("N9").Value = 77777
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = -30

# --- New row 10: Anita Ackermann-Abegger -----------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("A10").Value = 7560000000008
$ws.Range("B10").Value = "Ackermann-Abegger"
$ws.Range("C10").Value = "Anita"
$ws.Range("D10").Value = 34940
$ws.Range("E10").Value = "Sex#W"
$ws.Range("F10").Value = "Civil stat#B"
$ws.Range("G10").Value = "Gatte"
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = "CHResidenc#B"
$ws.Range("J10").Value = "Ackerstrasse"
$ws.Range("K10").Value = 11
$ws.Range("L10").Value = 804500
$ws.Range("M10").Value = "Zürich"
$ws.Range("N10").Value = 77777
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = -30

# --- Selection / active cell -----------------------------------------------------
[void]$ws.Range("C16").Select()
